$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set numeric-looking text values in column D (Price), preserving their
# original text (inline/shared string) cell type rather than becoming numbers.
function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "230.51"
Set-TextValue "D3" "22.44"
Set-TextValue "D4" "5.257"
Set-TextValue "D5" "0.05545"
Set-TextValue "D7" "6.477"
Set-TextValue "D8" "1.055"
Set-TextValue "D9" "0.7828"
Set-TextValue "D10" "0.1382"
Set-TextValue "D11" "0.07398"
Set-TextValue "D13" "0.02969"
Set-TextValue "D14" "0.09273"
Set-TextValue "D15" "0.001663"
Set-TextValue "D16" "3.257"
Set-TextValue "D17" "0.04792"
Set-TextValue "D18" "0.0005889"
Set-TextValue "D19" "0.006216"
Set-TextValue "D20" "0.005235"
Set-TextValue "D21" "0.001064"
Set-TextValue "D23" "3.917"
Set-TextValue "D24" "2.146"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
Set-TextValue "D40" "0.04004"
Set-TextValue "D41" "0.007083"
$ws.Range("E41").Value = "40KickTokenKICK"
Set-TextValue "D42" "0.003409"
Set-TextValue "D44" "0.009975"
Set-TextValue "D48" "0.04467"
